$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first three data rows (2007年, 2008年, 2009年); this shifts
# 2010年..2020年 up from rows 5-15 to rows 2-12.
$ws.Rows("2:4").Delete()

# Append the new 2021年 row. Copy the format of the row above (2020年, now
# row 12) into the new row 13 so the date cell keeps the bordered/bold/
# centered style used by the rest of column A.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$ws.Range("J13").Value = 1398
